$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''66.005.17'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -4.90%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.283.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -5.53%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.04%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''558.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -3.41%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''185.81'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -3.24%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +0.03%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  -2.82%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''3.276.91'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -5.34%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '''  -8.50%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.587'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -4.77%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''47.55'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -7.37%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  -6.81%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = '''BitcoinCash'
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = '''https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = '''636.06'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -2.66%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = '''Polkadot'
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = '''https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = '''8.63'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -5.43%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''3.814.38'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -5.42%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''66.014.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -4.74%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''18.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = '''  -3.32%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''3.286.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -5.61%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''11.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -7.86%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''0.907'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -3.90%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''18.37'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +1.75%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''107.77'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +8.85%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  -7.08%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  -7.33%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''  -7.03%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  -3.42%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  -6.43%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''30.36'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -6.57%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = '''  -6.04%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''6.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -6.72%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  -4.67%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  -3.64%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = '''OKB'
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = '''https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = '''57.68'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -5.18%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = '''Bittensor'
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = '''https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = '''526.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -0.03%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = '''Dai'
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = '''https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = '''1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -0.08%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = '''Maker'
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = '''https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = '''3.700.17'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -0.55%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = '''  -4.09%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '''  -7.76%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = '''Kaspa'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = '''0.130'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -2.49%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''2.71'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -6.74%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = '''InjectiveProtocol'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = '''33.02'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -3.60%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = '''CoreDAO'
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = '''https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = '''3.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -5.45%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '''  -9.14%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''3.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -1.59%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.0414'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -6.05%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  -3.64%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''2.60'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -7.89%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''1.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +0.02%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  +3.12%  '
$ws.Range("E51").Style = "Normal"
